# Invigilation Step 3 & 4 done.
#
# Only the cells whose displayed value actually changes are touched below.
# (A number of other <c> entries in the underlying XML get renumbered
# because the shared-string table shrinks/grows, but Excel/COM handles
# that bookkeeping automatically -- we just set the new literal values.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "     "
$ws.Range("E2").Value = " BIO543"

# Row 3
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "CSE344/CSE544/ ECE344/ECE544"
$ws.Range("G3").Value = "C2001,LHC"

# Row 4
$ws.Range("C4").Value = "      "

# Row 5
$ws.Range("F5").ClearContents()

# Row 6
$ws.Range("G6").ClearContents()

# Row 7
$ws.Range("B7").Value = "    "
$ws.Range("E7").Value = " MTH310 / MTH520      "

# Row 10
$ws.Range("A10").Value = 45414
$ws.Range("E10").Value = "   "

# Selection moved from E6 to E10
$ws.Range("E10").Select()
